$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellValue($row, $col, $oldText, $newText, $newBold) {
    $cell = $tbl.Cell($row, $col)
    $s = $cell.Range.Start
    $e = $cell.Range.End - 1
    $rng = $d.Range($s, $e)
    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 1) | Out-Null
    if ($null -ne $newBold) {
        $cell2 = $tbl.Cell($row, $col)
        if ($newBold) {
            $cell2.Range.Font.Bold = -1
        } else {
            $cell2.Range.Font.Bold = 0
        }
    }
}

Set-CellValue 2 2 "0.286" "0.115" $null
Set-CellValue 2 3 "0.630" "0.105" $null
Set-CellValue 6 3 "-0.209" "0.000" $false
Set-CellValue 7 2 "-0.049" "0.000" $false
Set-CellValue 8 2 "-0.097" "0.000" $false
Set-CellValue 8 3 "-0.132" "0.000" $false
Set-CellValue 9 2 "-0.266" "-0.183" $null
Set-CellValue 9 3 "-0.333" "0.000" $false
Set-CellValue 10 3 "0.032" "0.000" $false
Set-CellValue 12 2 "0.000" "-0.007" $true
Set-CellValue 13 3 "0.097" "0.000" $false
Set-CellValue 14 2 "-0.392" "-0.140" $null
Set-CellValue 14 3 "-0.489" "-0.080" $null
Set-CellValue 18 3 "0.038" "0.000" $false
Set-CellValue 19 3 "0.116" "0.000" $false
Set-CellValue 20 3 "0.116" "0.000" $false
Set-CellValue 21 3 "0.026" "0.000" $false
Set-CellValue 22 3 "0.026" "0.000" $false

Write-Host "Done"
